$d = $word.ActiveDocument

# Locate the "VON SPERLING..." bibliography paragraph; the three
# paragraphs that directly follow it (a blank spacer paragraph, the
# "Ver no Jupiter..." line, and the "(c) 2020 ..." footer line) are the
# ones that need to be removed, while the blank paragraph and the
# page-break paragraph that come after them must stay untouched.

$findRange = $d.Content
$found = $findRange.Find.Execute(
    "VON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos. 2. ed. Belo Horizonte: UFMG, 1996.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the VON SPERLING reference paragraph"
}

# Work out that paragraph's 1-based index within the document.
$prefix = $d.Range(0, $findRange.Start)
$anchorIndex = $prefix.Paragraphs.Count + 1

# Delete the three paragraphs that follow the anchor, starting from the
# highest index so earlier deletions don't shift later indices.
$d.Paragraphs.Item($anchorIndex + 3).Range.Delete()
$d.Paragraphs.Item($anchorIndex + 2).Range.Delete()
$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()
